$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.659958333333333
$ws.Range("H2").Value = 10.979875
$ws.Range("I2").Value = 0.4781132044744068
$ws.Range("J2").Value = 0.4781132044744067
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.679174
$ws.Range("N2").Value = 8.037521999999999
$ws.Range("O2").Value = 0.02942326717729479
$ws.Range("P2").Value = 0.02942326717729479
$ws.Range("Q2").Value = 9.80566520775
$ws.Range("R2").Value = 88.25098686974999
$ws.Range("S2").Value = 0.01406765255624305
$ws.Range("T2").Value = 0.01406765255624304
$ws.Range("G3").Value = 3.659958333333333
$ws.Range("H3").Value = 10.979875
$ws.Range("I3").Value = 0.4781132044744068
$ws.Range("J3").Value = 0.4781132044744067
$ws.Range("O3").Value = 0.2465847468531156
$ws.Range("P3").Value = 0.2465847468531155
$ws.Range("Q3").Value = 82.17739581433334
$ws.Range("R3").Value = 739.5965623289999
$ws.Range("S3").Value = 0.1178954234924535
$ws.Range("T3").Value = 0.1178954234924534
$ws.Range("G4").Value = 3.659958333333333
$ws.Range("H4").Value = 10.979875
$ws.Range("I4").Value = 0.4781132044744068
$ws.Range("J4").Value = 0.4781132044744067
$ws.Range("M4").Value = 65.67046766666668
$ws.Range("N4").Value = 197.011403
$ws.Range("O4").Value = 0.7212072511207682
$ws.Range("P4").Value = 0.7212072511207681
$ws.Range("Q4").Value = 240.351175390514
$ws.Range("R4").Value = 2163.160578514625
$ws.Range("S4").Value = 0.3448187099235287
$ws.Range("T4").Value = 0.3448187099235286
$ws.Range("G5").Value = 3.659958333333333
$ws.Range("H5").Value = 10.979875
$ws.Range("I5").Value = 0.4781132044744068
$ws.Range("J5").Value = 0.4781132044744067
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2535676666666667
$ws.Range("N5").Value = 0.760703
$ws.Range("O5").Value = 0.002784734848821526
$ws.Range("P5").Value = 0.002784734848821525
$ws.Range("Q5").Value = 0.9280470946805557
$ws.Range("R5").Value = 8.352423852125
$ws.Range("S5").Value = 0.001331418502181612
$ws.Range("T5").Value = 0.001331418502181612
$ws.Range("I6").Value = 0.3193330932870009
$ws.Range("J6").Value = 0.3193330932870008
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.679174
$ws.Range("N6").Value = 8.037521999999999
$ws.Range("O6").Value = 0.02942326717729479
$ws.Range("P6").Value = 0.02942326717729479
$ws.Range("Q6").Value = 6.549230126303999
$ws.Range("R6").Value = 58.94307113673599
$ws.Range("S6").Value = 0.009395822922335428
$ws.Range("T6").Value = 0.009395822922335428
$ws.Range("I7").Value = 0.3193330932870009
$ws.Range("J7").Value = 0.3193330932870008
$ws.Range("O7").Value = 0.2465847468531156
$ws.Range("P7").Value = 0.2465847468531155
$ws.Range("S7").Value = 0.07874266996999744
$ws.Range("T7").Value = 0.07874266996999743
$ws.Range("I8").Value = 0.3193330932870009
$ws.Range("J8").Value = 0.3193330932870008
$ws.Range("M8").Value = 65.67046766666668
$ws.Range("N8").Value = 197.011403
$ws.Range("O8").Value = 0.7212072511207682
$ws.Range("P8").Value = 0.7212072511207681
$ws.Range("Q8").Value = 160.531195529296
$ws.Range("R8").Value = 1444.780759763664
$ws.Range("S8").Value = 0.2303053424014097
$ws.Range("T8").Value = 0.2303053424014097
$ws.Range("I9").Value = 0.3193330932870009
$ws.Range("J9").Value = 0.3193330932870008
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.2535676666666667
$ws.Range("N9").Value = 0.760703
$ws.Range("O9").Value = 0.002784734848821526
$ws.Range("P9").Value = 0.002784734848821525
$ws.Range("Q9").Value = 0.619845146896
$ws.Range("R9").Value = 5.578606322064
$ws.Range("S9").Value = 0.0008892579932582866
$ws.Range("T9").Value = 0.0008892579932582862
$ws.Range("G10").Value = 1.388093333333333
$ws.Range("H10").Value = 4.16428
$ws.Range("I10").Value = 0.1813315046964271
$ws.Range("J10").Value = 0.1813315046964271
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.679174
$ws.Range("N10").Value = 8.037521999999999
$ws.Range("O10").Value = 0.02942326717729479
$ws.Range("P10").Value = 0.02942326717729479
$ws.Range("Q10").Value = 3.718943568239999
$ws.Range("R10").Value = 33.47049211416
$ws.Range("S10").Value = 0.005335365310343861
$ws.Range("T10").Value = 0.00533536531034386
$ws.Range("G11").Value = 1.388093333333333
$ws.Range("H11").Value = 4.16428
$ws.Range("I11").Value = 0.1813315046964271
$ws.Range("J11").Value = 0.1813315046964271
$ws.Range("O11").Value = 0.2465847468531156
$ws.Range("P11").Value = 0.2465847468531155
$ws.Range("Q11").Value = 31.16699287029333
$ws.Range("R11").Value = 280.5029358326399
$ws.Range("S11").Value = 0.04471358318206302
$ws.Range("T11").Value = 0.044713583182063
$ws.Range("G12").Value = 1.388093333333333
$ws.Range("H12").Value = 4.16428
$ws.Range("I12").Value = 0.1813315046964271
$ws.Range("J12").Value = 0.1813315046964271
$ws.Range("M12").Value = 65.67046766666668
$ws.Range("N12").Value = 197.011403
$ws.Range("O12").Value = 0.7212072511207682
$ws.Range("P12").Value = 0.7212072511207681
$ws.Range("Q12").Value = 91.15673836498223
$ws.Range("R12").Value = 820.41064528484
$ws.Range("S12").Value = 0.1307775960437029
$ws.Range("T12").Value = 0.1307775960437028
$ws.Range("G13").Value = 1.388093333333333
$ws.Range("H13").Value = 4.16428
$ws.Range("I13").Value = 0.1813315046964271
$ws.Range("J13").Value = 0.1813315046964271
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2535676666666667
$ws.Range("N13").Value = 0.760703
$ws.Range("O13").Value = 0.002784734848821526
$ws.Range("P13").Value = 0.002784734848821525
$ws.Range("Q13").Value = 0.3519755876488889
$ws.Range("R13").Value = 3.16778028884
$ws.Range("S13").Value = 0.0005049601603173847
$ws.Range("T13").Value = 0.0005049601603173846
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.162456
$ws.Range("H14").Value = 0.487368
$ws.Range("I14").Value = 0.02122219754216535
$ws.Range("J14").Value = 0.02122219754216534
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.679174
$ws.Range("N14").Value = 8.037521999999999
$ws.Range("O14").Value = 0.02942326717729479
$ws.Range("P14").Value = 0.02942326717729479
$ws.Range("Q14").Value = 0.435247891344
$ws.Range("R14").Value = 3.917231022096
$ws.Range("S14").Value = 0.0006244263883724598
$ws.Range("T14").Value = 0.0006244263883724597
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.162456
$ws.Range("H15").Value = 0.487368
$ws.Range("I15").Value = 0.02122219754216535
$ws.Range("J15").Value = 0.02122219754216534
$ws.Range("O15").Value = 0.2465847468531156
$ws.Range("P15").Value = 0.2465847468531155
$ws.Range("Q15").Value = 3.647640163776
$ws.Range("R15").Value = 32.828761473984
$ws.Range("S15").Value = 0.005233070208601654
$ws.Range("T15").Value = 0.005233070208601652
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.162456
$ws.Range("H16").Value = 0.487368
$ws.Range("I16").Value = 0.02122219754216535
$ws.Range("J16").Value = 0.02122219754216534
$ws.Range("M16").Value = 65.67046766666668
$ws.Range("N16").Value = 197.011403
$ws.Range("O16").Value = 0.7212072511207682
$ws.Range("P16").Value = 0.7212072511207681
$ws.Range("Q16").Value = 10.668561495256
$ws.Range("R16").Value = 96.01705345730402
$ws.Range("S16").Value = 0.01530560275212699
$ws.Range("T16").Value = 0.01530560275212699
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.162456
$ws.Range("H17").Value = 0.487368
$ws.Range("I17").Value = 0.02122219754216535
$ws.Range("J17").Value = 0.02122219754216534
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.2535676666666667
$ws.Range("N17").Value = 0.760703
$ws.Range("O17").Value = 0.002784734848821526
$ws.Range("P17").Value = 0.002784734848821525
$ws.Range("Q17").Value = 0.04119358885600001
$ws.Range("R17").Value = 0.370742299704
$ws.Range("S17").Value = 0.0000590981930642423705191
$ws.Range("T17").Value = 0.0000590981930642423501903
